# Update online job ads data to July 2025 data
#
# Row 13 of Sheet1 ("Job adverts by occupation") tracks the Textkernel /
# ONS online job-adverts series. The data has moved on a month, so the
# "Latest period (release date)" column (C) now shows the period that used
# to be shown as "Next period (release date)" (D), and column D is rolled
# forward to the new upcoming period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C13 ("Latest period") picks up what used to be the "next" period (D13),
# and loses the one-off small "Open Sans" formatting it had - it now uses
# the sheet's normal/default cell style, same as the rest of column C.
$ws.Range("C13").Value = $ws.Range("D13").Value2
$ws.Range("C13").Style = "Normal"

# D13 ("Next period") is rolled forward to the new upcoming release.
$ws.Range("D13").Value = "Aug 2025 (September (TBC))"
